$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.376.26'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '3.377.24'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.69%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '3.376.05'
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -0.40%  '
$ws.Range("E10").Value = '  -1.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.124'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("E12").Value = '  -1.18%  '
$ws.Range("D13").Value = '3.952.93'
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("E14").Value = '  +2.24%  '
$ws.Range("E15").Value = '  +0.88%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.05'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.62%  '
$ws.Range("D17").Value = '3.377.99'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '61.488.55'
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = '  +0.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '376.16'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.554'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.70%  '
$ws.Range("D24").Value = '3.513.75'
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("E26").Value = '  +4.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '71.55'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.99%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.79%  '
$ws.Range("E29").Value = '  -3.75%  '
$ws.Range("E30").Value = '  +0.13%  '
$ws.Range("E31").Value = '  +1.07%  '
$ws.Range("E32").Value = '  +2.79%  '
$ws.Range("E33").Value = '  +0.73%  '
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.58'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.27'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.61%  '
$ws.Range("E38").Value = '  -1.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '165.63'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0771'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.47%  '
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.72'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.773'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.71%  '
$ws.Range("E44").Value = '  +0.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.40'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.64%  '
$ws.Range("E46").Value = '  +0.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.39'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.17%  '
$ws.Range("E48").Value = '  -1.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.62'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.31%  '
$ws.Range("D50").Value = '2.361.48'
$ws.Range("E50").Value = '  +2.92%  '
$ws.Range("E51").Value = '  -0.76%  '
